# Append a newly scraped Lancers job listing (2025-12-23 18:28:04 JST run).
# The scraper re-writes the whole "取得日時" (fetched-at) column with the new
# run timestamp, inserts the freshly found listing right after the existing
# "AI" cluster (i.e. as the new row 6), and pushes the previously-existing
# rows 6-9 down to rows 7-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-12-23 18:28:04"

# --- Snapshot the current (pre-edit) rows 6-9 so we can shift them down to
#     rows 7-10 after the new row is written into row 6. (Columns A and E
#     don't need snapshotting: A gets the new run timestamp everywhere, and
#     the whole table only ever contains "期限情報なし" in column E.) ---
$oldB6 = $ws.Cells.Item(6, 2).Value2
$oldC6 = $ws.Cells.Item(6, 3).Value2
$oldD6 = $ws.Cells.Item(6, 4).Value2
$oldF6 = $ws.Cells.Item(6, 6).Value2
$oldG6 = $ws.Cells.Item(6, 7).Value2
$oldH6 = $ws.Cells.Item(6, 8).Value2

$oldB7 = $ws.Cells.Item(7, 2).Value2
$oldC7 = $ws.Cells.Item(7, 3).Value2
$oldD7 = $ws.Cells.Item(7, 4).Value2
$oldF7 = $ws.Cells.Item(7, 6).Value2
$oldG7 = $ws.Cells.Item(7, 7).Value2
$oldH7 = $ws.Cells.Item(7, 8).Value2

$oldB8 = $ws.Cells.Item(8, 2).Value2
$oldC8 = $ws.Cells.Item(8, 3).Value2
$oldD8 = $ws.Cells.Item(8, 4).Value2
$oldF8 = $ws.Cells.Item(8, 6).Value2
$oldG8 = $ws.Cells.Item(8, 7).Value2
$oldH8 = $ws.Cells.Item(8, 8).Value2

$oldB9 = $ws.Cells.Item(9, 2).Value2
$oldC9 = $ws.Cells.Item(9, 3).Value2
$oldD9 = $ws.Cells.Item(9, 4).Value2
$oldF9 = $ws.Cells.Item(9, 6).Value2
$oldG9 = $ws.Cells.Item(9, 7).Value2

# --- Refresh the "取得日時" column for every existing row (2-9) to the new
#     scrape run timestamp. ---
$ws.Range("A2:A9").Value = $timestamp

# --- Write the newly-discovered listing into row 6 (new entry). ---
$ws.Cells.Item(6, 2).Value = "初回 【AWSクラウドリフト】業務アプリ移行支援エンジニア募集(Java / .NET)"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5459847"
$ws.Cells.Item(6, 7).Value = 103
$ws.Cells.Item(6, 8).Value = "★Java ◇アプリ"

# --- Shift the previous rows 6-9 down into rows 7-10. ---
$ws.Cells.Item(7, 2).Value = $oldB6
$ws.Cells.Item(7, 3).Value = $oldC6
$ws.Cells.Item(7, 4).Value = $oldD6
$ws.Cells.Item(7, 6).Value = $oldF6
$ws.Cells.Item(7, 7).Value = $oldG6
$ws.Cells.Item(7, 8).Value = $oldH6

$ws.Cells.Item(8, 2).Value = $oldB7
$ws.Cells.Item(8, 3).Value = $oldC7
$ws.Cells.Item(8, 4).Value = $oldD7
$ws.Cells.Item(8, 6).Value = $oldF7
$ws.Cells.Item(8, 7).Value = $oldG7
$ws.Cells.Item(8, 8).Value = $oldH7

$ws.Cells.Item(9, 2).Value = $oldB8
$ws.Cells.Item(9, 3).Value = $oldC8
$ws.Cells.Item(9, 4).Value = $oldD8
$ws.Cells.Item(9, 6).Value = $oldF8
$ws.Cells.Item(9, 7).Value = $oldG8
$ws.Cells.Item(9, 8).Value = $oldH8

$ws.Cells.Item(10, 1).Value = $timestamp
$ws.Cells.Item(10, 2).Value = $oldB9
$ws.Cells.Item(10, 3).Value = $oldC9
$ws.Cells.Item(10, 4).Value = $oldD9
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = $oldF9
$ws.Cells.Item(10, 7).Value = $oldG9
$ws.Cells.Item(10, 8).ClearContents()

# --- Rebuild the hyperlinks collection so it matches the new row layout.
#     (The engine only supports wholesale add/remove, not in-place edits of
#     an existing hyperlink's target, so the cleanest way to get the right
#     F2..F10 -> URL mapping is to clear everything and re-add in order.) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5455098")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445159")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5445154")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5459299")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5459847")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5459200")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5459128")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5458992")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5459456")
